$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.222.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "'1.856.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'241.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").Value = "'0.6997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.76%  "

$ws.Range("D7").Value = "'0.9997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.07786"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.28%  "

$ws.Range("D9").Value = "'0.3078"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.19%  "

$ws.Range("D10").Value = "'23.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.13%  "

$ws.Range("D11").Value = "'0.07799"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.57%  "

$ws.Range("D12").Value = "'1.861.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").Value = "'5.111"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.43%  "

$ws.Range("D14").Value = "'92.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").Value = "'0.6890"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.24%  "

$ws.Range("D16").Value = "'6.552"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.21%  "

$ws.Range("D17").Value = "'0.000008481"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.39%  "

$ws.Range("D18").Value = "'29.215.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("D19").Value = "'248.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.96%  "

$ws.Range("D20").Value = "'2.107.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("E21").Value = "  -2.32%  "

$ws.Range("D22").Value = "'0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").Value = "'7.550"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.16%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "'0.1509"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.57%  "

$ws.Range("D26").Value = "'161.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.03%  "

$ws.Range("D27").Value = "'8.864"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.51%  "

$ws.Range("D28").Value = "'18.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.04%  "

$ws.Range("D29").Value = "'1.551"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.70%  "

$ws.Range("D30").Value = "'4.261"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.13%  "

$ws.Range("D31").Value = "'4.213"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.39%  "

$ws.Range("D32").Value = "'1.195"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.08%  "

$ws.Range("D33").Value = "'0.05233"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.98%  "

$ws.Range("D34").Value = "'0.7641"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.82%  "

$ws.Range("D35").Value = "'1.850"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.11%  "

$ws.Range("D36").Value = "'1.172"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.10%  "

$ws.Range("D37").Value = "'2.706"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("D38").Value = "'0.01862"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.47%  "

$ws.Range("D39").Value = "'1.228.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.31%  "

$ws.Range("D40").Value = "'2.729"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.26%  "

$ws.Range("D41").Value = "'0.9009"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").Value = "'109.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.01%  "

$ws.Range("D43").Value = "'0.9993"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").Value = "'5.544"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.53%  "

$ws.Range("D45").Value = "'2.006.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("E46").Value = "  -2.98%  "

$ws.Range("D47").Value = "'65.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.75%  "

$ws.Range("D48").Value = "'9.572"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.80%  "

$ws.Range("D49").Value = "'0.5183"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D50").Value = "'1.751"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.94%  "

$ws.Range("D51").Value = "'7.057"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.35%  "

